# "Generate Report for Handoff"
#
# 506cd9fa-dcd2-4c54-9f3a-d707b9167a78 has finished translation and its
# handoff package was (re)generated, so its status flips from
# "In Translation" to "Ready for handoff" with fresh handoff timestamps
# and new handoff-file names; df502825-399b-49a4-a02d-44e84c1c2dbc stays
# "In Translation". The two rows also swap display order/hyperlink text
# on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"  (A1:D3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md", "", "", "df502825-399b-49a4-a02d-44e84c1c2dbc.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/df502825-399b-49a4-a02d-44e84c1c2dbc.md", "", "", "506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md")

$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "2016-12-11 14:12:04"

$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-12-11 14:12:47"

# ---------------------------------------------------------------------
# Sheet "zh-cn"  (A1:K3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md", "", "", "df502825-399b-49a4-a02d-44e84c1c2dbc.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85b49b8d3f344ca7ffbbe4f543f8a2dd7512c21b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.ab51b88cb1258ae9ff1cf529282e2890aef81f11.zh-cn.xlf", "", "", "df502825-399b-49a4-a02d-44e84c1c2dbc.262b2258c0bd889407dcb108604d47cdbf34b3e3.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/df502825-399b-49a4-a02d-44e84c1c2dbc.md", "", "", "506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/df502825-399b-49a4-a02d-44e84c1c2dbc.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85b49b8d3f344ca7ffbbe4f543f8a2dd7512c21b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/df502825-399b-49a4-a02d-44e84c1c2dbc.262b2258c0bd889407dcb108604d47cdbf34b3e3.zh-cn.xlf", "", "", "506cd9fa-dcd2-4c54-9f3a-d707b9167a78.ab51b88cb1258ae9ff1cf529282e2890aef81f11.zh-cn.xlf")

$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "2016-03-11 14:11:57"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-11 14:12:44"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"  (A1:K3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md", "", "", "df502825-399b-49a4-a02d-44e84c1c2dbc.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4131ff8350daddff02bcc11ea0fadb33a3ae9814/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/506cd9fa-dcd2-4c54-9f3a-d707b9167a78.ab51b88cb1258ae9ff1cf529282e2890aef81f11.de-de.xlf", "", "", "df502825-399b-49a4-a02d-44e84c1c2dbc.262b2258c0bd889407dcb108604d47cdbf34b3e3.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/df502825-399b-49a4-a02d-44e84c1c2dbc.md", "", "", "506cd9fa-dcd2-4c54-9f3a-d707b9167a78.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/916926cda7dce7f0bcaf79d8250d5dc4e36755d3/e2e/df502825-399b-49a4-a02d-44e84c1c2dbc.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4131ff8350daddff02bcc11ea0fadb33a3ae9814/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/df502825-399b-49a4-a02d-44e84c1c2dbc.262b2258c0bd889407dcb108604d47cdbf34b3e3.de-de.xlf", "", "", "506cd9fa-dcd2-4c54-9f3a-d707b9167a78.ab51b88cb1258ae9ff1cf529282e2890aef81f11.de-de.xlf")

$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "2016-03-11 14:12:04"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-11 14:12:47"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# Restore original active sheet / selection.
$wb.Worksheets.Item("Overview").Activate()
